$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01134666666666667
$ws.Range("H2").Value = 0.03404
$ws.Range("I2").Value = 0.001209510404472147
$ws.Range("J2").Value = 0.001209510404472147
$ws.Range("M2").Value = 0.04647766666666667
$ws.Range("N2").Value = 0.139433
$ws.Range("O2").Value = 0.1662164546338858
$ws.Range("P2").Value = 0.1662164546338858
$ws.Range("Q2").Value = 0.0005273665911111112
$ws.Range("R2").Value = 0.00474629932
$ws.Range("S2").Value = 0.0002010405312741575
$ws.Range("T2").Value = 0.0002010405312741575
# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01134666666666667
$ws.Range("H3").Value = 0.03404
$ws.Range("I3").Value = 0.001209510404472147
$ws.Range("J3").Value = 0.001209510404472147
$ws.Range("O3").Value = 0.4511935188540693
$ws.Range("P3").Value = 0.4511935188540693
$ws.Range("Q3").Value = 0.001431533288888889
$ws.Range("R3").Value = 0.0128837996
$ws.Range("S3").Value = 0.0005457232554843965
$ws.Range("T3").Value = 0.0005457232554843965
# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01134666666666667
$ws.Range("H4").Value = 0.03404
$ws.Range("I4").Value = 0.001209510404472147
$ws.Range("J4").Value = 0.001209510404472147
$ws.Range("M4").Value = 0.1069803333333333
$ws.Range("N4").Value = 0.320941
$ws.Range("O4").Value = 0.3825900265120449
$ws.Range("P4").Value = 0.3825900265120449
$ws.Range("Q4").Value = 0.001213870182222222
$ws.Range("R4").Value = 0.01092483164
$ws.Range("S4").Value = 0.0004627466177135928
$ws.Range("T4").Value = 0.0004627466177135928
# Row 5
$ws.Range("I5").Value = 0.8865539289740954
$ws.Range("J5").Value = 0.8865539289740952
$ws.Range("M5").Value = 0.04647766666666667
$ws.Range("N5").Value = 0.139433
$ws.Range("O5").Value = 0.1662164546338858
$ws.Range("P5").Value = 0.1662164546338858
$ws.Range("Q5").Value = 0.3865522128875556
$ws.Range("R5").Value = 3.478969915988
$ws.Range("S5").Value = 0.1473598509158159
$ws.Range("T5").Value = 0.1473598509158159
# Row 6
$ws.Range("I6").Value = 0.8865539289740954
$ws.Range("J6").Value = 0.8865539289740952
$ws.Range("O6").Value = 0.4511935188540693
$ws.Range("P6").Value = 0.4511935188540693
$ws.Range("S6").Value = 0.4000073868677227
$ws.Range("T6").Value = 0.4000073868677226
# Row 7
$ws.Range("I7").Value = 0.8865539289740954
$ws.Range("J7").Value = 0.8865539289740952
$ws.Range("M7").Value = 0.1069803333333333
$ws.Range("N7").Value = 0.320941
$ws.Range("O7").Value = 0.3825900265120449
$ws.Range("P7").Value = 0.3825900265120449
$ws.Range("Q7").Value = 0.8897495840751114
$ws.Range("R7").Value = 8.007746256676002
$ws.Range("S7").Value = 0.3391866911905567
$ws.Range("T7").Value = 0.3391866911905567
# Row 8
$ws.Range("G8").Value = 1.052914333333334
$ws.Range("H8").Value = 3.158743
$ws.Range("I8").Value = 0.1122365606214325
$ws.Range("J8").Value = 0.1122365606214325
$ws.Range("M8").Value = 0.04647766666666667
$ws.Range("N8").Value = 0.139433
$ws.Range("O8").Value = 0.1662164546338858
$ws.Range("P8").Value = 0.1662164546338858
$ws.Range("Q8").Value = 0.04893700141322223
$ws.Range("R8").Value = 0.440433012719
$ws.Range("S8").Value = 0.01865556318679571
$ws.Range("T8").Value = 0.01865556318679571
# Row 9
$ws.Range("G9").Value = 1.052914333333334
$ws.Range("H9").Value = 3.158743
$ws.Range("I9").Value = 0.1122365606214325
$ws.Range("J9").Value = 0.1122365606214325
$ws.Range("O9").Value = 0.4511935188540693
$ws.Range("P9").Value = 0.4511935188540693
$ws.Range("Q9").Value = 0.1328391820077778
$ws.Range("R9").Value = 1.19555263807
$ws.Range("S9").Value = 0.0506404087308622
$ws.Range("T9").Value = 0.0506404087308622
# Row 10
$ws.Range("G10").Value = 1.052914333333334
$ws.Range("H10").Value = 3.158743
$ws.Range("I10").Value = 0.1122365606214325
$ws.Range("J10").Value = 0.1122365606214325
$ws.Range("M10").Value = 0.1069803333333333
$ws.Range("N10").Value = 0.320941
$ws.Range("O10").Value = 0.3825900265120449
$ws.Range("P10").Value = 0.3825900265120449
$ws.Range("Q10").Value = 0.1126411263514445
$ws.Range("R10").Value = 1.013770137163
$ws.Range("S10").Value = 0.0429405887037746
$ws.Range("T10").Value = 0.0429405887037746
